$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 245
$ws1.Cells.Item(7, 6).Value = 6451
$ws1.Cells.Item(10, 6).Value = 129
$ws1.Cells.Item(11, 6).Value = 5758
$ws1.Cells.Item(12, 6).Value = 35
$ws1.Cells.Item(14, 6).Value = 1227
$ws1.Cells.Item(15, 6).Value = 1227
$ws1.Cells.Item(19, 6).Value = 82
$ws1.Cells.Item(21, 6).Value = 328
$ws1.Cells.Item(24, 6).Value = 4100
$ws1.Cells.Item(25, 6).Value = 26

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 245
$ws4.Cells.Item(7, 6).Value = 6452
$ws4.Cells.Item(10, 6).Value = 129
$ws4.Cells.Item(11, 6).Value = 5758
$ws4.Cells.Item(12, 6).Value = 35
$ws4.Cells.Item(14, 6).Value = 1227
$ws4.Cells.Item(15, 6).Value = 1227
$ws4.Cells.Item(19, 6).Value = 82
$ws4.Cells.Item(21, 6).Value = 328
$ws4.Cells.Item(24, 6).Value = 4100
$ws4.Cells.Item(26, 6).Value = 26
